$p = $ppt.ActivePresentation

# Ensure notes master exists
$nm = $p.NotesMaster

$s1 = $p.Slides.Item(1)
$s2 = $p.Slides.Item(2)

# Slide 1: empty notes slide (just touch it so it gets created)
$np1 = $s1.NotesPage
$np1.Shapes.Placeholders.Item(2).TextFrame.TextRange.Text = ""

# Slide 2: notes text
$np2 = $s2.NotesPage
$np2.Shapes.Placeholders.Item(2).TextFrame.TextRange.Text = "H1 Slide 2`nH2 bla`nH3 blub"
